$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bus voltage magnitude results (case with 380 kV) for rows 2-25
# Columns: B,C,D,E,F,I,J,K,L,M,N (G/H unchanged)
$data = @(
    @(1.02, 1.054601782356865, 1.062033843852216, 0.992614727750844, 1.070775148022582, 1.05360258256022, 1.059613021632177, 1.064756287437982, 0.9955398523335997, 1.073474087112851, 1.023593548019571),
    @(1.02, 1.05563299783454, 1.062879824408843, 0.9936372048519299, 1.071751971969555, 1.053951911540173, 1.060295187620437, 1.065416926584839, 0.9963617723202687, 1.07426694082838, 1.023826566979097),
    @(1.02, 1.056300231080739, 1.063427182863084, 0.9942998659930998, 1.072384313232432, 1.054176657681791, 1.060735948988808, 1.065843710015767, 0.9968940712668347, 1.074779618983051, 1.023976994951966),
    @(1.02, 1.056580728247904, 1.063657280422334, 0.994578699834602, 1.072650213822048, 1.05427083120998, 1.060921090239376, 1.066022963218707, 0.9971179600053012, 1.074995064387412, 1.024040150709913),
    @(1.02, 1.056627824505698, 1.063695914099007, 0.994625531979634, 1.072694863457893, 1.054286625190469, 1.06095216718322, 1.066053050843676, 0.9971555583673455, 1.075031233652262, 1.024050749899597),
    @(1.02, 1.056303979127817, 1.063430257485522, 0.994303590798249, 1.072387865957028, 1.054177917250419, 1.060738423463952, 1.065846105860696, 0.9968970624462089, 1.074782498107595, 1.023977839173193),
    @(1.02, 1.054950293107242, 1.062319756451521, 0.9929600610674297, 1.07110521385545, 1.053720907550299, 1.059843696349347, 1.064979697079691, 0.9958175282591056, 1.073742108213444, 1.023672370421873),
    @(1.02, 1.05256467882557, 1.060362571790902, 0.9906006454969559, 1.068847113892528, 1.052905702748027, 1.058262139175852, 1.063447671518091, 0.9939188001724441, 1.071906136122054, 1.023131415196243),
    @(1.02, 1.050974099663244, 1.059057583810289, 0.989033133672735, 1.067343151726689, 1.052355589289683, 1.057204463532115, 1.062422773810849, 0.9926553831429383, 1.070680379805182, 1.022768985130365),
    @(1.02, 1.050285319240481, 1.058492468301749, 0.988355674866747, 1.066692265648369, 1.052115809853371, 1.056745696623038, 1.061978143510053, 0.9921088820399291, 1.070149196654328, 1.022611625208434),
    @(1.02, 1.05002946786342, 1.058282552539153, 0.9881042295826724, 1.066450548898807, 1.052026508269612, 1.056575172094612, 1.061812861556071, 0.9919059725120875, 1.069951828313198, 1.022553110928812),
    @(1.02, 1.050084349203354, 1.058327580480023, 0.9881581567098651, 1.066502395634617, 1.05204567447843, 1.056611755526652, 1.061848320816374, 0.9919494934313052, 1.06999416735184, 1.022565665325813),
    @(1.02, 1.050264170643687, 1.058475116733856, 0.9883348863814464, 1.066672284217668, 1.052108432987441, 1.056731603420529, 1.061964483831042, 0.9920921077337197, 1.070132883415055, 1.022606789700708),
    @(1.02, 1.05037496354469, 1.058566017812432, 0.9884438009545853, 1.06677696499077, 1.052147069196153, 1.056805430019947, 1.062036038945528, 0.9921799884222134, 1.07021834258276, 1.022632119334485),
    @(1.02, 1.051019810692331, 1.059095087705132, 0.9890781214508737, 1.067386356071505, 1.052371469446881, 1.057234893819079, 1.06245226471597, 0.9926916645766087, 1.070715623797014, 1.022779419637839),
    @(1.02, 1.051424292917479, 1.059426946926197, 0.989476357848556, 1.067768702106928, 1.052511807711203, 1.057504074629865, 1.062713126610541, 0.9930127773699352, 1.071027442183961, 1.022871703429087),
    @(1.02, 1.051660215832348, 1.059620510225545, 0.9897087662937556, 1.06799175081257, 1.052593512516189, 1.057661007363905, 1.062865201608419, 0.9932001317071769, 1.071209279786286, 1.022925489960065),
    @(1.02, 1.05174065865671, 1.059686509521953, 0.9897880325774034, 1.068067810199258, 1.052621345937716, 1.057714504498166, 1.062917041452493, 0.9932640239640975, 1.071271274784669, 1.022943822812488),
    @(1.02, 1.051380896293167, 1.059391342042325, 0.9894336180360679, 1.067727676607784, 1.052496766478695, 1.057475201924251, 1.062685147009296, 0.9929783193494215, 1.0709939912392, 1.022861806496187),
    @(1.02, 1.050211217925104, 1.058431671166318, 0.9882828385668249, 1.066622254845054, 1.052089958707846, 1.056696314461689, 1.061930280208205, 0.9920501090198102, 1.07009203673659, 1.022594681357636),
    @(1.02, 1.049475750256923, 1.05782824941288, 0.9875604150241495, 1.065927529392357, 1.051832812347514, 1.056205913909358, 1.061454933893999, 0.9914670000341481, 1.06952457620211, 1.022426360114225),
    @(1.02, 1.049865639888636, 1.058148138233782, 0.9879432794643023, 1.066295788232654, 1.051969260356555, 1.056465949245666, 1.061706993261516, 0.991776070289318, 1.069825432381801, 1.022515625339829),
    @(1.02, 1.051400505374117, 1.059407430370818, 0.9894529299347244, 1.067746214163717, 1.052503563434214, 1.057488248491463, 1.062697790037953, 0.9929938892766442, 1.071009106407204, 1.022866278620953),
    @(1.02, 1.053181446921732, 1.060868588090712, 0.9912096547607049, 1.06943063547857, 1.053117624403429, 1.058671592948174, 1.06384436339314, 0.9944092447426414, 1.072381093343518, 1.023271581806308)
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}
